$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 22:30"

# Update country case statistics (values refreshed + a few rows re-ranked by total cases,
# which swaps the country name shown in column A for the affected row pairs).

# Row 4
$ws.Range("B4").Value = 3273514
$ws.Range("C4").Value = 53515
$ws.Range("D4").Value = 1443423
$ws.Range("E4").Value = 1693660
$ws.Range("G4").Value = 609
$ws.Range("H4").Value = 136431

# Row 6
$ws.Range("B6").Value = 822570
$ws.Range("C6").Value = 27728
$ws.Range("E6").Value = 284220

# Row 19
$ws.Range("B19").Value = 199475
$ws.Range("C19").Value = 277
$ws.Range("E19").Value = 6346
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 9129

# Row 27
$ws.Range("B27").Value = 80235
$ws.Range("C27").Value = 981
$ws.Range("D27").Value = 23274
$ws.Range("E27").Value = 53259
$ws.Range("G27").Value = 85
$ws.Range("H27").Value = 3702

# Row 50
$ws.Range("E50").Value = 4904
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 104

# Row 70
$ws.Range("A70").Value = "Costa de Marfil"
$ws.Range("B70").Value = 12052
$ws.Range("C70").Value = 302
$ws.Range("D70").Value = 6080
$ws.Range("E70").Value = 5891
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 81

# Row 71
$ws.Range("A71").Value = "Uzbekistan"
$ws.Range("B71").Value = 11857
$ws.Range("C71").Value = 293
$ws.Range("D71").Value = 7335
$ws.Range("E71").Value = 4468
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 54

# Row 85
$ws.Range("A85").Value = "Costa Rica"
$ws.Range("B85").Value = 6845
$ws.Range("C85").Value = 360
$ws.Range("D85").Value = 2110
$ws.Range("E85").Value = 4709
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 26

# Row 86
$ws.Range("A86").Value = "Bulgaria"
$ws.Range("B86").Value = 6672
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 3229
$ws.Range("E86").Value = 3181
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 262

# Row 87
$ws.Range("A87").Value = "Haiti"
$ws.Range("B87").Value = 6582
$ws.Range("C87").Value = 96
$ws.Range("D87").Value = 2459
$ws.Range("E87").Value = 3993
$ws.Range("G87").Value = 7
$ws.Range("H87").Value = 130

# Row 90
$ws.Range("B90").Value = 5969
$ws.Range("C90").Value = 88
$ws.Range("D90").Value = 4732
$ws.Range("E90").Value = 1200
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 37

# Row 95
$ws.Range("B95").Value = 4968
$ws.Range("C95").Value = 13
$ws.Range("D95").Value = 4689
$ws.Range("E95").Value = 223

# Row 97
$ws.Range("A97").Value = "Republica de Africa Central"
$ws.Range("B97").Value = 4259
$ws.Range("C97").Value = 59
$ws.Range("D97").Value = 1142
$ws.Range("E97").Value = 3064
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 53

# Row 98
$ws.Range("A98").Value = "Hungria"
$ws.Range("B98").Value = 4223
$ws.Range("C98").Value = 3
$ws.Range("D98").Value = 2941
$ws.Range("E98").Value = 689
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 593

# Row 132
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 1252
$ws.Range("C132").Value = 42
$ws.Range("D132").Value = 635
$ws.Range("E132").Value = 614
$ws.Range("H132").Value = 3

# Row 133
$ws.Range("A133").Value = "Tunez"
$ws.Range("B133").Value = 1240
$ws.Range("C133").Value = 9
$ws.Range("D133").Value = 1067
$ws.Range("E133").Value = 123
$ws.Range("H133").Value = 50

# Row 145
$ws.Range("B145").Value = 942
$ws.Range("C145").Value = 16
$ws.Range("D145").Value = 320
$ws.Range("E145").Value = 609
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 13

# Row 150
$ws.Range("A150").Value = "Surinam"
$ws.Range("B150").Value = 726
$ws.Range("C150").Value = 32
$ws.Range("D150").Value = 468
$ws.Range("E150").Value = 241
$ws.Range("H150").Value = 17

# Row 151
$ws.Range("A151").Value = "Crucero"
$ws.Range("B151").Value = 712
$ws.Range("D151").Value = 651
$ws.Range("E151").Value = 48
$ws.Range("H151").Value = 13

# Row 152
$ws.Range("A152").Value = "Togo"
$ws.Range("B152").Value = 710
$ws.Range("C152").Value = 6
$ws.Range("D152").Value = 494
$ws.Range("E152").Value = 201
$ws.Range("H152").Value = 15

# Row 208
$ws.Range("B208").Value = 17
$ws.Range("C208").Value = 1
$ws.Range("E208").Value = 2
